# Crackteck-modules.xlsx — "added auth and leads apis"
#
# Updates the status tracker sheet: several rows move from
# Pending/Not Done/Some-Changes-pending to Done, one row moves from
# Not Done to Some-Changes-pending, row 31 gets a taller custom row
# height, a stray description string in C80 is corrected back to the
# section title, and the view's last selection/scroll position is
# updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status column (E) updates -------------------------------------------
# Style 65 = "Done" (green fill), Style 67 = "Some Changes are pending"
# (amber fill), Style 68 = "Pending"/"Not Done" (red fill/font). There is
# no named cell style for these, so we copy the formatting from an
# existing cell that already carries the desired style, then overwrite
# the text.

# E35: Pending -> Done
$ws.Range("E7").Copy()
$ws.Range("E35").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E35").Value = "Done"

# E36: Some Changes are pending -> Done
$ws.Range("E7").Copy()
$ws.Range("E36").PasteSpecial(-4122)
$ws.Range("E36").Value = "Done"

# E59: Not Done -> Some Changes are pending
$ws.Range("E23").Copy()
$ws.Range("E59").PasteSpecial(-4122)
$ws.Range("E59").Value = "Some Changes are pending"

# E65: Some Changes are pending -> Done
$ws.Range("E7").Copy()
$ws.Range("E65").PasteSpecial(-4122)
$ws.Range("E65").Value = "Done"

# E66: Not Done -> Done
$ws.Range("E7").Copy()
$ws.Range("E66").PasteSpecial(-4122)
$ws.Range("E66").Value = "Done"

$excel.CutCopyMode = $false

# --- Row 31: taller custom row height -------------------------------------
$ws.Rows.Item(31).RowHeight = 20.25

# --- C80: fix stray value back to the section title -----------------------
$ws.Range("C80").Value = "Low Stock Reports"

# --- View state: scroll position + last selected cell ---------------------
[void]$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 52
$win.ScrollColumn = 1
[void]$ws.Range("E59").Select()
